$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 18. This shifts the existing rows 18-41 down
# to rows 19-42 (each keeps its own data), and leaves row 18 ready to be
# filled in with the new weekly record.
$ws.Rows.Item(18).Insert()

# The row Insert already copied formatting from the row above (so D18
# keeps the date style). Now fill in the static (repeated) columns for
# the new record, matching the rest of this product's rows.
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = "Terminal La Palmera de La Serena"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44880
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100101
$ws.Range("H18").Value = "Berries"
$ws.Range("I18").Value = 100101001
$ws.Range("J18").Value = "Arándano (blue)"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 7000
$ws.Range("O18").Value = 8000
$ws.Range("P18").Value = 7500
$ws.Range("Q18").Value = "$/bandeja 2 kilos"
$ws.Range("R18").Value = "Provincia de Limarí"
$ws.Range("S18").Value = 3750
$ws.Range("T18").Value = 2
